$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.881.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'1.907.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'313.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.14%  "
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").Value = "'0.5005"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.77%  "
$ws.Range("D8").Value = "'0.3815"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "'0.07280"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.16%  "
$ws.Range("D10").Value = "'0.9084"
$ws.Range("D10").Style = "Normal"
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("D12").Value = "'0.07663"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.86%  "
$ws.Range("D13").Value = "'1.914.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("D14").Value = "'5.480"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.77%  "
$ws.Range("D15").Value = "'91.67"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("D16").Value = "'1.004"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").Value = "'0.000008725"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.19%  "
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("D19").Value = "'27.906.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.95%  "
$ws.Range("E20").Value = "  -2.30%  "
$ws.Range("D21").Value = "'5.165"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("D23").Value = "'6.603"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.52%  "
$ws.Range("D24").Value = "'153.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.39%  "
$ws.Range("D25").Value = "'1.882"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.88%  "
$ws.Range("D26").Value = "'2.238"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.98%  "
$ws.Range("D27").Value = "'18.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.93%  "
$ws.Range("E28").Value = "  -0.90%  "
$ws.Range("D29").Value = "'4.908"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.17%  "
$ws.Range("D30").Value = "'0.08970"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").Value = "'3.216"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.78%  "
$ws.Range("D32").Value = "'1.234"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.61%  "
$ws.Range("D33").Value = "'0.7663"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.05%  "
$ws.Range("D34").Value = "'4.646"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.91%  "
$ws.Range("D35").Value = "'0.02063"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("D36").Value = "'2.557"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.44%  "
$ws.Range("D37").Value = "'0.5619"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.14%  "
$ws.Range("E38").Value = "  -0.93%  "
$ws.Range("D39").Value = "'3.012"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.80%  "
$ws.Range("D40").Value = "'0.05251"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.31%  "
$ws.Range("D41").Value = "'6.963"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.71%  "
$ws.Range("D42").Value = "'8.495"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").Value = "'0.1516"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.91%  "
$ws.Range("D44").Value = "'111.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.26%  "
$ws.Range("E45").Value = "  -1.18%  "
$ws.Range("D46").Value = "'0.4810"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.51%  "
$ws.Range("D47").Value = "'1.003"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.21%  "
$ws.Range("E48").Value = "  -1.46%  "
$ws.Range("D50").Value = "'0.06071"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.19%  "
$ws.Range("D51").Value = "'0.9003"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.16%  "
